# Updated cryptos list on Mon Nov  6 09:52:13 UTC 2023 with GitHub Actions
#
# The "Price" (D) and "Volume(1h)" (E) columns are refreshed with new
# scraped values. Two coin pairs also swapped rank/position, which moves
# their Coin name (B) and Link (C) between the two rows in addition to
# their Price/Volume.
#
# All of these columns hold plain TEXT in the workbook (prices are written
# with "." as a thousands separator, e.g. "35.462.53", so they can't be
# numeric cells) - so for any new value that *looks* like a genuine number
# we force the cell to Text format first, otherwise Excel's COM layer would
# silently reinterpret it as a numeric/date value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value, taken from the diff.
$updates = [ordered]@{
    "D2"  = "35.462.53"
    "E2"  = "  -0.09%  "
    "D3"  = "1.906.75"
    "E3"  = "  +0.01%  "
    "E4"  = "  +0.04%  "
    "D5"  = "0.704"
    "E5"  = "  +10.43%  "
    "D6"  = "246.66"
    "E7"  = "  +0.05%  "
    "D8"  = "40.74"
    "E8"  = "  -3.34%  "
    "E9"  = "  +4.08%  "
    "D10" = "52.70"
    "E10" = "  +8.19%  "
    "D11" = "0.0727"
    "E11" = "  +3.13%  "
    "E12" = "  -0.79%  "
    "D13" = "2.184.22"
    "E13" = "  +0.10%  "
    "D14" = "12.55"
    "E14" = "  +1.72%  "
    "D15" = "0.714"
    "E15" = "  +2.18%  "
    "D16" = "1.906.21"
    "E16" = "  +0.84%  "
    "E17" = "  +1.44%  "
    "D18" = "35.456.38"
    "E18" = "  -0.22%  "
    "D19" = "73.11"
    "E19" = "  +1.64%  "
    "E20" = "  -0.20%  "
    "D21" = "242.18"
    "E21" = "  -0.53%  "
    "E22" = "  +1.18%  "
    "D23" = "5.05"
    "E23" = "  +3.97%  "
    "E24" = "  -0.03%  "
    "E25" = "  +1.10%  "
    "E26" = "  +5.20%  "
    "D27" = "169.21"
    "E27" = "  -1.82%  "
    "E28" = "  +0.76%  "
    "D29" = "18.93"
    "E29" = "  +5.23%  "
    "E30" = "  +5.08%  "
    "E32" = "  +2.80%  "
    "D33" = "0.0574"
    "E33" = "  +0.54%  "
    "D34" = "4.20"
    "E34" = "  +0.78%  "

    # Row 35/36: BinanceUSD and WEMIXToken swapped places.
    "B35" = "BinanceUSD"
    "C35" = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
    "D35" = "1.01"
    "E35" = "  +0.02%  "
    "B36" = "WEMIXToken"
    "C36" = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
    "D36" = "1.87"
    "E36" = "  +5.85%  "

    "D37" = "0.923"
    "E37" = "  -5.89%  "
    "D38" = "1.47"
    "E38" = "  +9.35%  "
    "D39" = "2.01"
    "E39" = "  -0.95%  "
    "D40" = "96.22"
    "E40" = "  +5.58%  "
    "E41" = "  +0.75%  "
    "D42" = "16.57"
    "E42" = "  +5.18%  "

    # Row 43/44: Kaspa and VeChain swapped places.
    "B43" = "VeChain"
    "C43" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "D43" = "0.0209"
    "E43" = "  +2.02%  "
    "B44" = "Kaspa"
    "C44" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "D44" = "0.0649"
    "E44" = "  +1.27%  "

    "D45" = "1.353.51"
    "E45" = "  +0.23%  "
    "E46" = "  +0.52%  "
    "D47" = "46.36"
    "E47" = "  -8.05%  "
    "D48" = "2.42"
    "E48" = "  +0.25%  "
    "E49" = "  +1.18%  "
    "D50" = "12.18"
    "E50" = "  -5.13%  "
    "E51" = "  -1.86%  "
}

foreach ($ref in $updates.Keys) {
    $value = $updates[$ref]
    $cell = $ws.Range($ref)

    # Prices like "246.66" or "0.0727" round-trip as genuine numbers through
    # COM automation, which would flip the cell from text to numeric and
    # reformat/round the displayed value. Anything else (the two-separator
    # "35.462.53"-style prices, the percentage strings, and the coin/link
    # names) already fails that numeric parse and stays text on its own.
    if ($value -match '^[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = "@"
    }

    $cell.Value = $value
}
